# "Chiffres COVID-19 Valais" daily data-entry update.
# Source-data edits only; all dependent formulas (cumulative columns
# B, H, J, K) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Disable concurrent (multi-threaded) calculation, matching the
# workbook-level calcPr concurrentCalc="0" setting in the saved file.
try {
    $excel.MultiThreadedCalculation.Enabled = $false
} catch {
}

# Helper: columns L and M are formatted as Text ("@"), so a plain
# .Value assignment gets stored as a text string instead of a real
# number. Flip the format to General for the write, then restore it,
# so the stored cell keeps its original Text number format but holds
# a genuine numeric value (as Excel itself does when typing digits
# into a Text-formatted cell that was previously numeric).
function Set-NumericValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# --- Corrections to already-entered daily rows ---
$ws.Range("C244").Value = 882
$ws.Range("C245").Value = 887
$ws.Range("C275").Value = 142

# --- Update of "Patients COVID-19 hospitalisés hors SI" (col G) ---
$ws.Range("G340").Value = 104
$ws.Range("G341").Value = 113
$ws.Range("G342").Value = 117

# --- Row 343 (2020-08-16) revisions ---
$ws.Range("C343").Value = 101
$ws.Range("G343").Value = 122
Set-NumericValue $ws.Range("L343") 2
Set-NumericValue $ws.Range("M343") 4

# --- Row 344 (2020-08-17) revisions ---
$ws.Range("C344").Value = 84
$ws.Range("G344").Value = 119
Set-NumericValue $ws.Range("L344") 1

# --- Row 345 (2020-08-18) newly entered daily figures ---
$ws.Range("C345").Value = 7
$ws.Range("E345").Value = 12
$ws.Range("F345").Value = 7
$ws.Range("G345").Value = 115
Set-NumericValue $ws.Range("L345") 0
Set-NumericValue $ws.Range("M345") 0

# --- Active-cell selection moved to P7 on the frozen bottom-right pane ---
$ws.Range("P7").Select() | Out-Null
